$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 265.0247
$ws.Range("J17").Value = 265.0247
$ws.Range("L17").Value = 795.0741
$ws.Range("N17").Value = -1131.0741
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H39").Value = 144.78572
$ws.Range("I39").Value = 182.66667
$ws.Range("J39").Value = 76.59999999999999
$ws.Range("K39").Value = 548.00001
$ws.Range("L39").Value = 229.8
$ws.Range("M39").Value = -252.00001
$ws.Range("N39").Value = -821.8
$ws.Range("H40").Value = 1196.9318
$ws.Range("I40").Value = 1089.619
$ws.Range("J40").Value = 1294.9131
$ws.Range("K40").Value = 1089.619
$ws.Range("L40").Value = 1294.9131
$ws.Range("M40").Value = -914.6189999999999
$ws.Range("N40").Value = -1644.9131
$ws.Range("H51").Value = 3498.3333
$ws.Range("I51").Value = 3225.5334
$ws.Range("J51").Value = 3953
$ws.Range("K51").Value = 3225.5334
$ws.Range("L51").Value = 3953
$ws.Range("M51").Value = -2741.5334
$ws.Range("N51").Value = -4921
$ws.Range("H116").Value = 40701.715
$ws.Range("I116").Value = 64862.35
$ws.Range("J116").Value = 3362.5454
$ws.Range("K116").Value = 64862.35
$ws.Range("L116").Value = 3362.5454
$ws.Range("M116").Value = -61420.35
$ws.Range("N116").Value = -10246.5454

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 1448.35
$ws.Range("I2").Value = 1313.0667
$ws.Range("J2").Value = 1854.2
$ws.Range("K2").Value = 1313.0667
$ws.Range("L2").Value = 1854.2
$ws.Range("M2").Value = -1200.0667
$ws.Range("N2").Value = -2080.2
$ws.Range("H31").Value = 6990.3335
$ws.Range("I31").Value = 6990.3335
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6990.3335
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -6696.3335
$ws.Range("N31").ClearContents()
$ws.Range("H32").Value = 961374.6
$ws.Range("I32").Value = 1228681.2
$ws.Range("J32").Value = 4698.1055
$ws.Range("K32").Value = 1228681.2
$ws.Range("L32").Value = 4698.1055
$ws.Range("M32").Value = -1228394.2
$ws.Range("N32").Value = -5272.1055
$ws.Range("H63").Value = 4263.75
$ws.Range("I63").Value = 4263.75
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 4263.75
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3577.75
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 4263.75
$ws.Range("I66").Value = 4263.75
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 21318.75
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -17886.75
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 213191.25
$ws.Range("I74").Value = 239227.72
$ws.Range("K74").Value = 239227.72
$ws.Range("M74").Value = -238353.72
$ws.Range("H77").Value = 213191.25
$ws.Range("I77").Value = 239227.72
$ws.Range("K77").Value = 1196138.6
$ws.Range("M77").Value = -1191770.6
$ws.Range("H101").Value = 39467.25
$ws.Range("J101").Value = 39467.25
$ws.Range("L101").Value = 39467.25
$ws.Range("N101").Value = -45957.25
$ws.Range("H116").Value = 1448.35
$ws.Range("I116").Value = 1313.0667
$ws.Range("J116").Value = 1854.2
$ws.Range("K116").Value = 1313.0667
$ws.Range("L116").Value = 1854.2
$ws.Range("M116").Value = 980.9332999999999
$ws.Range("N116").Value = -6442.2

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 1448.35
$ws.Range("I3").Value = 1313.0667
$ws.Range("J3").Value = 1854.2
$ws.Range("K3").Value = 1313.0667
$ws.Range("L3").Value = 1854.2
$ws.Range("M3").Value = -1199.0667
$ws.Range("N3").Value = -2082.2
$ws.Range("H82").Value = 11279.667
$ws.Range("I82").Value = 6575.6
$ws.Range("J82").Value = 34800
$ws.Range("K82").Value = 6575.6
$ws.Range("L82").Value = 34800
$ws.Range("M82").Value = -6192.6
$ws.Range("N82").Value = -35566
$ws.Range("H85").Value = 11279.667
$ws.Range("I85").Value = 6575.6
$ws.Range("J85").Value = 34800
$ws.Range("K85").Value = 6575.6
$ws.Range("L85").Value = 34800
$ws.Range("M85").Value = -5249.6
$ws.Range("N85").Value = -37452
$ws.Range("H94").Value = 720.1799999999999
$ws.Range("I94").Value = 638.5278
$ws.Range("J94").Value = 930.1429000000001
$ws.Range("K94").Value = 638.5278
$ws.Range("L94").Value = 930.1429000000001
$ws.Range("M94").Value = -187.5278
$ws.Range("N94").Value = -1832.1429

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 2282.7576
$ws.Range("I31").Value = 1667.5128
$ws.Range("J31").Value = 3171.4443
$ws.Range("K31").Value = 1667.5128
$ws.Range("L31").Value = 3171.4443
$ws.Range("M31").Value = -1372.5128
$ws.Range("N31").Value = -3761.4443
$ws.Range("H34").Value = 2282.7576
$ws.Range("I34").Value = 1667.5128
$ws.Range("J34").Value = 3171.4443
$ws.Range("K34").Value = 1667.5128
$ws.Range("L34").Value = 3171.4443
$ws.Range("M34").Value = -1465.5128
$ws.Range("N34").Value = -3575.4443
$ws.Range("H99").Value = 68962.87
$ws.Range("I99").Value = 144522.86
$ws.Range("J99").Value = 2847.875
$ws.Range("K99").Value = 144522.86
$ws.Range("L99").Value = 2847.875
$ws.Range("M99").Value = -143024.86
$ws.Range("N99").Value = -5843.875
$ws.Range("H126").Value = 68962.87
$ws.Range("I126").Value = 144522.86
$ws.Range("J126").Value = 2847.875
$ws.Range("K126").Value = 433568.58
$ws.Range("L126").Value = 8543.625
$ws.Range("M126").Value = -431098.58
$ws.Range("N126").Value = -13483.625
$ws.Range("H132").Value = 1971.0426
$ws.Range("I132").Value = 1178.4828
$ws.Range("J132").Value = 3247.9443
$ws.Range("K132").Value = 3535.4484
$ws.Range("L132").Value = 9743.832900000001
$ws.Range("M132").Value = -1005.4484
$ws.Range("N132").Value = -14803.8329
$ws.Range("H135").Value = 57688.25
$ws.Range("J135").Value = 40592.668
$ws.Range("L135").Value = 40592.668
$ws.Range("N135").Value = -50732.668

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H13").Value = 1400
$ws.Range("I13").Value = 798
$ws.Range("J13").Value = 2002
$ws.Range("K13").Value = 2394
$ws.Range("L13").Value = 6006
$ws.Range("M13").Value = -2226
$ws.Range("N13").Value = -6342
$ws.Range("H113").Value = 8929106
$ws.Range("I113").Value = 13158417
$ws.Range("J113").Value = 561
$ws.Range("K113").Value = 39475251
$ws.Range("L113").Value = 1683
$ws.Range("M113").Value = -39473081
$ws.Range("N113").Value = -6023
$ws.Range("H137").Value = 5728.5312
$ws.Range("I137").Value = 1764
$ws.Range("J137").Value = 9226.647000000001
$ws.Range("K137").Value = 5292
$ws.Range("L137").Value = 27679.941
$ws.Range("M137").Value = -192
$ws.Range("N137").Value = -37879.94100000001

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 3908.963
$ws.Range("I102").Value = 2642.8572
$ws.Range("J102").Value = 8340.333000000001
$ws.Range("K102").Value = 2642.8572
$ws.Range("L102").Value = 8340.333000000001
$ws.Range("M102").Value = -1020.8572
$ws.Range("N102").Value = -11584.333
$ws.Range("H122").Value = 1279.6897
$ws.Range("I122").Value = 1247.7059
$ws.Range("J122").Value = 1325
$ws.Range("K122").Value = 3743.1177
$ws.Range("L122").Value = 3975
$ws.Range("M122").Value = -1293.1177
$ws.Range("N122").Value = -8875
$ws.Range("H126").Value = 3257.5386
$ws.Range("I126").Value = 3062.2222
$ws.Range("J126").Value = 3360.9412
$ws.Range("K126").Value = 9186.6666
$ws.Range("L126").Value = 10082.8236
$ws.Range("M126").Value = -6716.6666
$ws.Range("N126").Value = -15022.8236

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H55").Value = 199.55
$ws.Range("I55").Value = 175
$ws.Range("J55").Value = 236.375
$ws.Range("K55").Value = 175
$ws.Range("L55").Value = 236.375
$ws.Range("M55").Value = -2
$ws.Range("N55").Value = -582.375
$ws.Range("H132").Value = 6809.4443
$ws.Range("I132").Value = 2549.9473
$ws.Range("J132").Value = 9922.154
$ws.Range("K132").Value = 7649.841899999999
$ws.Range("L132").Value = 29766.462
$ws.Range("M132").Value = -5119.841899999999
$ws.Range("N132").Value = -34826.462
